$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“머신러닝 파워드 애플리케이션”이 곧 출간될 예정입니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/08/25/%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-%ed%8c%8c%ec%9b%8c%eb%93%9c-%ec%95%a0%ed%94%8c%eb%a6%ac%ec%bc%80%ec%9d%b4%ec%85%98%ec%9d%b4-%ea%b3%a7-%ec%b6%9c%ea%b0%84%eb%90%a0-%ec%98%88%ec%a0%95%ec%9e%85/"

$ws.Range("D16").Value = "Combinational Class Activation Maps for Weakly Supervised Object Localization 내용 정리 [XAI-18]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/167"

$ws.Range("D37").Value = "[Paper Review] Anomaly Detection of Time Series with Smoothness-Inducing Sequential Variational Auto-Encoder"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1815&mod=document&pageid=1"

$ws.Range("D46").Value = "[유한양행] 2021년 08월, 생물정보학(Bioinformatics 채용), 합성신약 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/414"
